# Auto-generated Excel COM-interop script to apply Golem_Profits data refresh
# Updates currentAveragePrice/LevePrice/LeveProfit columns (H-N) across all 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 307.3
$ws.Range("I6").Value = 307.3
$ws.Range("K6").Value = 921.9000000000001
$ws.Range("M6").Value = -809.9000000000001
$ws.Range("H33").Value = 77.5
$ws.Range("I33").Value = 77.5
$ws.Range("K33").Value = 77.5
$ws.Range("M33").Value = 151.5
$ws.Range("H42").Value = 68.5
$ws.Range("I42").Value = 68.5
$ws.Range("K42").Value = 205.5
$ws.Range("M42").Value = 24.5
$ws.Range("H55").Value = 413.07693
$ws.Range("I55").Value = 523.55554
$ws.Range("J55").Value = 164.5
$ws.Range("K55").Value = 523.55554
$ws.Range("L55").Value = 164.5
$ws.Range("M55").Value = -309.55554
$ws.Range("N55").Value = -592.5
$ws.Range("H62").Value = 2998.3333
$ws.Range("I62").Value = 2998.3333
$ws.Range("K62").Value = 2998.3333
$ws.Range("M62").Value = -2374.3333
$ws.Range("H65").Value = 2998.3333
$ws.Range("I65").Value = 2998.3333
$ws.Range("K65").Value = 14991.6665
$ws.Range("M65").Value = -11871.6665
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2064
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10320
$ws.Range("H97").Value = 950
$ws.Range("J97").Value = 950
$ws.Range("L97").Value = 2850
$ws.Range("N97").Value = -3842
$ws.Range("H98").Value = 24999.428
$ws.Range("I98").Value = 20999.2
$ws.Range("K98").Value = 20999.2
$ws.Range("M98").Value = -19501.2
$ws.Range("H115").Value = 1509
$ws.Range("I115").Value = 1509
$ws.Range("K115").Value = 4527
$ws.Range("M115").Value = -2960
$ws.Range("H118").Value = 549
$ws.Range("I118").Value = 549
$ws.Range("K118").Value = 1647
$ws.Range("M118").Value = 10
$ws.Range("H122").Value = 24999.428
$ws.Range("I122").Value = 20999.2
$ws.Range("K122").Value = 62997.60000000001
$ws.Range("M122").Value = -60547.60000000001
$ws.Range("H127").Value = 1331.3334
$ws.Range("I127").Value = 1331.3334
$ws.Range("K127").Value = 3994.0002
$ws.Range("M127").Value = 965.9998000000001
$ws.Range("H137").Value = 2747.5
$ws.Range("I137").Value = 2495
$ws.Range("K137").Value = 7485
$ws.Range("M137").Value = -4935
$ws.Range("H138").Value = 1090.4
$ws.Range("I138").Value = 656
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 1968
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = 3172
$ws.Range("N138").Value = -25280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 14250
$ws.Range("I41").Value = 3000
$ws.Range("J41").Value = 18000
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 18000
$ws.Range("M41").Value = -2586
$ws.Range("N41").Value = -18828
$ws.Range("H88").Value = 3169.889
$ws.Range("J88").Value = 3478.625
$ws.Range("L88").Value = 3478.625
$ws.Range("N88").Value = -4290.625
$ws.Range("H91").Value = 3169.889
$ws.Range("J91").Value = 3478.625
$ws.Range("L91").Value = 3478.625
$ws.Range("N91").Value = -6286.625
$ws.Range("H132").Value = 2971
$ws.Range("I132").Value = 2442.5
$ws.Range("K132").Value = 7327.5
$ws.Range("M132").Value = -4797.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 10300.25
$ws.Range("I7").Value = 400.33334
$ws.Range("K7").Value = 400.33334
$ws.Range("M7").Value = -287.33334
$ws.Range("H107").Value = 1876.4375
$ws.Range("I107").Value = 1850.2307
$ws.Range("K107").Value = 1850.2307
$ws.Range("M107").Value = 69.76929999999993
$ws.Range("H134").Value = 2703
$ws.Range("I134").Value = 2823.6
$ws.Range("J134").Value = 2100
$ws.Range("K134").Value = 8470.799999999999
$ws.Range("L134").Value = 6300
$ws.Range("M134").Value = -5935.799999999999
$ws.Range("N134").Value = -11370

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 33334000
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 100000000
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 100000000
$ws.Range("M4").Value = -888
$ws.Range("N4").Value = -100000224
$ws.Range("H31").Value = 804.2727
$ws.Range("I31").Value = 804.2727
$ws.Range("K31").Value = 804.2727
$ws.Range("M31").Value = -509.2727
$ws.Range("H34").Value = 804.2727
$ws.Range("I34").Value = 804.2727
$ws.Range("K34").Value = 804.2727
$ws.Range("M34").Value = -602.2727
$ws.Range("H99").Value = 1501500
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1501500
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1501500
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -1504496
$ws.Range("H126").Value = 1501500
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1501500
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 4504500
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -4509440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 400
$ws.Range("I8").Value = 400
$ws.Range("K8").Value = 1200
$ws.Range("M8").Value = -1061
$ws.Range("H12").Value = 190.44444
$ws.Range("I12").Value = 49
$ws.Range("K12").Value = 147
$ws.Range("M12").Value = 26
$ws.Range("H32").Value = 27381.117
$ws.Range("I32").Value = 6855
$ws.Range("J32").Value = 41749.4
$ws.Range("K32").Value = 20565
$ws.Range("L32").Value = 125248.2
$ws.Range("M32").Value = -20282
$ws.Range("N32").Value = -125814.2
$ws.Range("H129").Value = 2693.3333
$ws.Range("J129").Value = 5011
$ws.Range("L129").Value = 15033
$ws.Range("N129").Value = -25033
$ws.Range("H132").Value = 99
$ws.Range("I132").Value = 99
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 891
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 1639
$ws.Range("N132").ClearContents()
$ws.Range("H139").Value = 2000
$ws.Range("I139").Value = 2000
$ws.Range("K139").Value = 6000
$ws.Range("M139").Value = -860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 51524.625
$ws.Range("I102").Value = 58656.715
$ws.Range("K102").Value = 58656.715
$ws.Range("M102").Value = -57034.715
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2198.5
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 366.66666
$ws.Range("I18").Value = 366.66666
$ws.Range("K18").Value = 366.66666
$ws.Range("M18").Value = -194.66666
$ws.Range("H68").Value = 2125
$ws.Range("I68").Value = 1970
$ws.Range("K68").Value = 1970
$ws.Range("M68").Value = -1221
$ws.Range("H71").Value = 2125
$ws.Range("I71").Value = 1970
$ws.Range("K71").Value = 9850
$ws.Range("M71").Value = -6106
$ws.Range("H132").Value = 2699.75
$ws.Range("I132").Value = 2699.75
$ws.Range("K132").Value = 8099.25
$ws.Range("M132").Value = -5569.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4000
$ws.Range("I81").Value = 4000
$ws.Range("K81").Value = 8000
$ws.Range("M81").Value = -6939
$ws.Range("H84").Value = 4000
$ws.Range("I84").Value = 4000
$ws.Range("K84").Value = 40000
$ws.Range("M84").Value = -34696
$ws.Range("H132").Value = 1099
$ws.Range("I132").Value = 1099
$ws.Range("K132").Value = 3297
$ws.Range("M132").Value = -767
$ws.Range("H136").Value = 15545
$ws.Range("I136").Value = 15545
$ws.Range("K136").Value = 46635
$ws.Range("M136").Value = -44085
